# Update results and utils
#
# Rewrites the per-error-bin result rows (3-8) on the active sheet:
#   - Column A: the bin label changes from a text range (e.g. "0-5") to the
#     plain numeric lower bound (e.g. 5).
#   - Columns C:J (delta_k.* metrics, excluding the GroundTruth column B
#     which stays 0): refreshed with new values.
#   - Columns K:S (delta_eps.* metrics): cleared out entirely.
#   - Columns T:AB (delta_cov.*) and column B are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (A, C, D, E, F, G, H, I, J)
$rows = @{
    3 = @(5,  -8,   -3,  -2,   -16,   -2.5, -3,    -12, -10)
    4 = @(10, -7.5, -4,  -4.5, -17,   -2,   -7.5,  -2,  -8.5)
    5 = @(15, -13.5,-9,  -10,  -21,   -5,   -10.5, -3,  -14)
    6 = @(20, -13,  -4,  -6,   -19,   0,    -16,   -12, -15)
    7 = @(25, -14.5,-7,  -9,   -16.5, -2,   -15.5, 2,   -15.5)
    8 = @(30, -20,  -5,  -13,  -20,   -1,   -19,   3,   -19)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
    $ws.Range("H$r").Value = $vals[6]
    $ws.Range("I$r").Value = $vals[7]
    $ws.Range("J$r").Value = $vals[8]
    $ws.Range("K${r}:S${r}").ClearContents()
}
